# Generate Report for Handback
# Updates the "zh-cn" and "de-de" status sheets for the row corresponding to
# 62548f44-df53-44d5-880d-8d72a43105de.md: a new handback was produced, so
# the "Latest Target File", "Latest Handback File", "Latest Handback DateTime"
# and "Error Detail" columns (I, J, K, P) get populated for row 7.

$wb = $excel.ActiveWorkbook

$notLatestMsg = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1657f104443e34e393ebb39f88f49542cc38b101/e2e/62548f44-df53-44d5-880d-8d72a43105de.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9f7a2bbc43b53f5551e2d583c70a6b787f57ed28/e2e/62548f44-df53-44d5-880d-8d72a43105de.md."

# ---------------------------------------------------------------------------
# zh-cn sheet, row 7
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("I7").Value = "62548f44-df53-44d5-880d-8d72a43105de.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/9f7a2bbc43b53f5551e2d583c70a6b787f57ed28/e2e/62548f44-df53-44d5-880d-8d72a43105de.md", "", "", "62548f44-df53-44d5-880d-8d72a43105de.md") | Out-Null

$wsZh.Range("J7").Value = "62548f44-df53-44d5-880d-8d72a43105de.43e9d80c1113cb6a3bb4a0a206cef4a0539b858d.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-17 12:52:57"
$wsZh.Range("P7").Value = $notLatestMsg

# ---------------------------------------------------------------------------
# de-de sheet, row 7
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("I7").Value = "62548f44-df53-44d5-880d-8d72a43105de.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/9f7a2bbc43b53f5551e2d583c70a6b787f57ed28/e2e/62548f44-df53-44d5-880d-8d72a43105de.md", "", "", "62548f44-df53-44d5-880d-8d72a43105de.md") | Out-Null

$wsDe.Range("J7").Value = "62548f44-df53-44d5-880d-8d72a43105de.43e9d80c1113cb6a3bb4a0a206cef4a0539b858d.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-17 12:53:12"
$wsDe.Range("P7").Value = $notLatestMsg
